$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 63 to grow the table from 62 to 67 data rows
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()
$ws.Rows.Item(63).Insert()

# Pre-seed the new shared strings in the exact order the target workbook
# expects them appended to xl/sharedStrings.xml (matches commit diff order).
$ws.Range("C57").Value = "The tree has been adopted."
$ws.Range("C58").Value = "The tree has ben unadopted."
$ws.Range("C55").Value = "New tree has been created."
$ws.Range("C51").Value = "New note has been posted."
$ws.Range("C56").Value = "Press the button to delete this tree."
$ws.Range("C67").Value = "The tree has been deleted."

# Rewrite rows 50-67 (Note/Tree/Authentification/MySQL blocks) with the final data
$ws.Range("A50").Value = 604
$ws.Range("B50").Value = "Note"
$ws.Range("C50").Value = "The note has been updated."
$ws.Range("D50").Value = "Response"
$ws.Range("A51").Value = 605
$ws.Range("B51").Value = "Note"
$ws.Range("C51").Value = "New note has been posted."
$ws.Range("D51").Value = "Response"
$ws.Range("A52").Value = 606
$ws.Range("B52").Value = "Note"
$ws.Range("C52").Value = "Press the button to delete this note."
$ws.Range("D52").Value = "Message"
$ws.Range("A53").Value = 607
$ws.Range("B53").Value = "Note"
$ws.Range("C53").Value = "The note has been deleted."
$ws.Range("D53").Value = "Response"
$ws.Range("A54").Value = 634
$ws.Range("B54").Value = "Tree"
$ws.Range("C54").Value = "The tree has been updated."
$ws.Range("D54").Value = "Response"
$ws.Range("A55").Value = 635
$ws.Range("B55").Value = "Tree"
$ws.Range("C55").Value = "New tree has been created."
$ws.Range("D55").Value = "Response"
$ws.Range("A56").Value = 636
$ws.Range("B56").Value = "Tree"
$ws.Range("C56").Value = "Press the button to delete this tree."
$ws.Range("D56").Value = "Message"
$ws.Range("A57").Value = 638
$ws.Range("B57").Value = "Tree"
$ws.Range("C57").Value = "The tree has been adopted."
$ws.Range("D57").Value = "Response"
$ws.Range("A58").Value = 639
$ws.Range("B58").Value = "Tree"
$ws.Range("C58").Value = "The tree has ben unadopted."
$ws.Range("D58").Value = "Response"
$ws.Range("A59").Value = 900
$ws.Range("B59").Value = "Authentification"
$ws.Range("C59").Value = "Not logged in."
$ws.Range("D59").Value = "Response"
$ws.Range("A60").Value = 901
$ws.Range("B60").Value = "Authentification"
$ws.Range("C60").Value = "Access is not authorized."
$ws.Range("D60").Value = "Response"
$ws.Range("A61").Value = 902
$ws.Range("B61").Value = "Authentification"
$ws.Range("C61").Value = "Failed to log in."
$ws.Range("D61").Value = "Response"
$ws.Range("A62").Value = 903
$ws.Range("B62").Value = "Authentification"
$ws.Range("C62").Value = "Account is already exist."
$ws.Range("D62").Value = "Response"
$ws.Range("A63").Value = 904
$ws.Range("B63").Value = "Authentification"
$ws.Range("C63").Value = "Not a valid e-mail address."
$ws.Range("D63").Value = "Response"
$ws.Range("A64").Value = 905
$ws.Range("B64").Value = "Authentification"
$ws.Range("C64").Value = "New parent has been registered."
$ws.Range("D64").Value = "Response"
$ws.Range("A65").Value = 906
$ws.Range("B65").Value = "Authentification"
$ws.Range("C65").Value = "Not an admin account."
$ws.Range("D65").Value = "Response"
$ws.Range("A66").Value = "42S22"
$ws.Range("B66").Value = "MySQL"
$ws.Range("C66").Value = "Unknown column in a database table."
$ws.Range("D66").Value = "Response"
$ws.Range("A67").Value = 637
$ws.Range("B67").Value = "Tree"
$ws.Range("C67").Value = "The tree has been deleted."
$ws.Range("D67").Value = "Response"

# Update the active selection to match the target view state
$ws.Range("F64").Select()
